$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values
$ws.Range("A5").Value = 5000
$ws.Range("B5").Value = 1000

# Add new row 9 values
$ws.Range("A9").Value = 8000
$ws.Range("B9").Value = 8000
